$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 2 (Hall 2 and Hall 5 entries, Hall 1 value also changes)
$ws.Range("B2").Value = "MEC0811"
$ws.Range("C2").Value = "CIE2802"
$ws.Range("D2").Value = "GEN0801"

# Fill in new rows 3-7 with test data
$data = @(
    @("9AM - 12PM", "GEN0810", "GEN0807", "CIE3801"),
    @("9AM - 12PM", "GEN1801", "GEN0806", "CIE3804"),
    @("9AM - 12PM", "CIE1803", "GEN1805", "CIE4818"),
    @("9AM - 12PM", "GEN2810", "GEN1809", "GEN0809"),
    @("9AM - 12PM", "GEN0802", "CIE1808", "POW1804")
)

$row = 3
foreach ($rowData in $data) {
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $row++
}
